$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 81.07692
$ws.Range("I2").Value = 41.454544
$ws.Range("J2").Value = 299
$ws.Range("K2").Value = 41.454544
$ws.Range("L2").Value = 299
$ws.Range("M2").Value = 71.545456
$ws.Range("N2").Value = -525
$ws.Range("H4").Value = 362.83334
$ws.Range("I4").Value = 335.4
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 335.4
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -221.4
$ws.Range("N4").Value = -728
$ws.Range("H48").Value = 1272.2222
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 2225
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 6675
$ws.Range("M48").Value = -2708
$ws.Range("N48").Value = -7259
$ws.Range("H56").Value = 1272.2222
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 2225
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 6675
$ws.Range("M56").Value = -2466
$ws.Range("N56").Value = -7743
$ws.Range("H132").Value = 8005.089
$ws.Range("I132").Value = 5327.846
$ws.Range("K132").Value = 15983.538
$ws.Range("M132").Value = -13453.538
$ws.Range("H137").Value = 3730.5405
$ws.Range("I137").Value = 3797.8823
$ws.Range("K137").Value = 11393.6469
$ws.Range("M137").Value = -8843.6469
$ws.Range("H138").Value = 2674.0942
$ws.Range("I138").Value = 2608.3076
$ws.Range("K138").Value = 7824.9228
$ws.Range("M138").Value = -2684.9228

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 7000
$ws.Range("I26").Value = 7000
$ws.Range("K26").Value = 7000
$ws.Range("M26").Value = -6670
$ws.Range("H52").Value = 12500
$ws.Range("I52").Value = 12500
$ws.Range("K52").Value = 12500
$ws.Range("M52").Value = -12182
$ws.Range("H55").Value = 32999.332
$ws.Range("J55").Value = 32999.332
$ws.Range("L55").Value = 32999.332
$ws.Range("N55").Value = -33629.332
$ws.Range("H74").Value = 2335.077
$ws.Range("I74").Value = 2335.077
$ws.Range("K74").Value = 2335.077
$ws.Range("M74").Value = -1461.077
$ws.Range("H77").Value = 2335.077
$ws.Range("I77").Value = 2335.077
$ws.Range("K77").Value = 11675.385
$ws.Range("M77").Value = -7307.385000000002
$ws.Range("H132").Value = 3750.6155
$ws.Range("I132").Value = 3750.6155
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11251.8465
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8721.8465

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3207.111
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3406
$ws.Range("H122").Value = 13604.223
$ws.Range("I122").Value = 15117.25
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 45351.75
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -42901.75
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 7719.3716
$ws.Range("I132").Value = 1658.6957
$ws.Range("K132").Value = 4976.0871
$ws.Range("M132").Value = -2446.0871
$ws.Range("H134").Value = 3937.9375
$ws.Range("I134").Value = 3858
$ws.Range("K134").Value = 11574
$ws.Range("M134").Value = -9039
$ws.Range("H136").Value = 3207.111
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 545.8
$ws.Range("I4").Value = 549.0714
$ws.Range("K4").Value = 1647.2142
$ws.Range("M4").Value = -1535.2142

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15887429
$ws.Range("I11").Value = 5303001
$ws.Range("K11").Value = 5303001
$ws.Range("M11").Value = -5302862
$ws.Range("H12").Value = 9800
$ws.Range("J12").Value = 9800
$ws.Range("L12").Value = 9800
$ws.Range("N12").Value = -10080
$ws.Range("H29").Value = 8249.25
$ws.Range("I29").Value = 4500
$ws.Range("J29").Value = 9499
$ws.Range("K29").Value = 4500
$ws.Range("L29").Value = 9499
$ws.Range("M29").Value = -4210
$ws.Range("N29").Value = -10079
$ws.Range("H126").Value = 4371.1763
$ws.Range("J126").Value = 4287.4
$ws.Range("L126").Value = 12862.2
$ws.Range("N126").Value = -17802.2
$ws.Range("H132").Value = 3209.3914
$ws.Range("I132").Value = 2216.7693
$ws.Range("J132").Value = 4499.8
$ws.Range("K132").Value = 6650.3079
$ws.Range("L132").Value = 13499.4
$ws.Range("M132").Value = -4120.3079
$ws.Range("N132").Value = -18559.4
$ws.Range("H137").Value = 75000
$ws.Range("I137").Value = 75000
$ws.Range("K137").Value = 75000
$ws.Range("M137").Value = -69900

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 945.2857
$ws.Range("I22").Value = 943.3333
$ws.Range("J22").Value = 946.75
$ws.Range("K22").Value = 943.3333
$ws.Range("L22").Value = 946.75
$ws.Range("M22").Value = -648.3333
$ws.Range("N22").Value = -1536.75
$ws.Range("H26").Value = 20369.666
$ws.Range("I26").Value = 15554.5
$ws.Range("K26").Value = 15554.5
$ws.Range("M26").Value = -15259.5
$ws.Range("H27").Value = 945.2857
$ws.Range("I27").Value = 943.3333
$ws.Range("J27").Value = 946.75
$ws.Range("K27").Value = 943.3333
$ws.Range("L27").Value = 946.75
$ws.Range("M27").Value = -836.3333
$ws.Range("N27").Value = -1160.75
$ws.Range("H34").Value = 13000
$ws.Range("J34").Value = 13000
$ws.Range("L34").Value = 13000
$ws.Range("N34").Value = -13344
$ws.Range("H133").Value = 109703.164
$ws.Range("J133").Value = 109703.164
$ws.Range("L133").Value = 109703.164
$ws.Range("N133").Value = -114763.164

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20470
$ws.Range("H29").Value = 71500
$ws.Range("I29").Value = 45000
$ws.Range("K29").Value = 45000
$ws.Range("M29").Value = -44710
$ws.Range("H35").Value = 20000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20580
$ws.Range("H132").Value = 2953.5103
$ws.Range("I132").Value = 2841.3865
$ws.Range("K132").Value = 8524.1595
$ws.Range("M132").Value = -5994.1595

Write-Host "All updates applied."